# fix(department): make edit works and create department import with standalone component
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused trailing columns F:G and the extra data row 4.
$ws.Columns("F:G").Delete()
$ws.Rows("4").Delete()

# Header row: rename comp_code -> company_code.
$ws.Range("A1").Value = "company_code"

# Row 2: update the date value (keep it text, like the original import data).
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "2024-10-21"
$ws.Range("C2").Style = "Normal"

# Row 3: now a holiday entry (was an "event" entry) with new branch/date/description.
$ws.Range("B3").Value = "B002"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "2024-10-22"
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").Value = "holiday"
$ws.Range("E3").Value = "hari libu"
